# Update cryptos list with latest price/volume snapshot
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text like "51.960.78" - force text format so Excel
# does not reinterpret the assigned strings as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "52.014.29"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").Value = "2.780.67"
$ws.Range("E3").Value = "  -1.52%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "356.30"
$ws.Range("E5").Value = "  +0.15%  "

$ws.Range("D6").Value = "109.19"
$ws.Range("E6").Value = "  -3.82%  "

$ws.Range("D7").Value = "0.565"
$ws.Range("E7").Value = "  +2.10%  "

$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  -1.43%  "

$ws.Range("D10").Value = "40.21"
$ws.Range("E10").Value = "  -3.97%  "

$ws.Range("D11").Value = "0.0852"
$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("E12").Value = "  +0.77%  "

$ws.Range("D13").Value = "19.41"
$ws.Range("E13").Value = "  -3.19%  "

$ws.Range("D14").Value = "7.61"
$ws.Range("E14").Value = "  -1.47%  "

$ws.Range("D15").Value = "3.212.61"
$ws.Range("E15").Value = "  -1.02%  "

$ws.Range("D16").Value = "2.788.31"
$ws.Range("E16").Value = "  -1.46%  "

$ws.Range("D17").Value = "0.930"
$ws.Range("E17").Value = "  +3.52%  "

$ws.Range("D18").Value = "51.844.27"
$ws.Range("E18").Value = "  -0.10%  "

$ws.Range("D19").Value = "7.42"
$ws.Range("E19").Value = "  -0.02%  "

$ws.Range("D20").Value = "3.14"
$ws.Range("E20").Value = "  -0.74%  "

$ws.Range("D21").Value = "13.06"
$ws.Range("E21").Value = "  -4.01%  "

$ws.Range("D22").Value = "0.0₃0976"
$ws.Range("E22").Value = "  -2.03%  "

$ws.Range("D23").Value = "274.57"
$ws.Range("E23").Value = "  +1.56%  "

$ws.Range("D24").Value = "69.89"
$ws.Range("E24").Value = "  +0.14%  "

$ws.Range("D25").Value = "2.73"
$ws.Range("E25").Value = "  -2.05%  "

$ws.Range("D26").Value = "26.59"
$ws.Range("E26").Value = "  -0.87%  "

$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").Value = "10.13"
$ws.Range("E28").Value = "  -1.89%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  -1.23%  "

$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "0.143"
$ws.Range("E30").Value = "  +2.12%  "

$ws.Range("D31").Value = "0.0467"
$ws.Range("E31").Value = "  +1.60%  "

$ws.Range("D32").Value = "51.63"
$ws.Range("E32").Value = "  +1.33%  "

$ws.Range("D33").Value = "33.91"
$ws.Range("E33").Value = "  -0.20%  "

$ws.Range("D34").Value = "5.71"
$ws.Range("E34").Value = "  -2.38%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.0843"
$ws.Range("E35").Value = "  +1.33%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "5.32"
$ws.Range("E36").Value = "  +8.76%  "

$ws.Range("D38").Value = "3.22"
$ws.Range("E38").Value = "  +0.30%  "

$ws.Range("D39").Value = "18.08"
$ws.Range("E39").Value = "  -1.73%  "

$ws.Range("E40").Value = "  -4.53%  "

$ws.Range("D41").Value = "2.54"
$ws.Range("E41").Value = "  -1.55%  "

$ws.Range("D42").Value = "0.115"
$ws.Range("E42").Value = "  -0.39%  "

$ws.Range("E43").Value = "  -2.96%  "

$ws.Range("D44").Value = "121.39"
$ws.Range("E44").Value = "  -5.40%  "

$ws.Range("D45").Value = "21.92"
$ws.Range("E45").Value = "  -7.26%  "

$ws.Range("D46").Value = "2.057.70"
$ws.Range("E46").Value = "  -0.97%  "

$ws.Range("D47").Value = "3.25"
$ws.Range("E47").Value = "  -2.86%  "

$ws.Range("E48").Value = "  -1.91%  "

$ws.Range("D49").Value = "5.70"
$ws.Range("E49").Value = "  +0.09%  "

$ws.Range("D50").Value = "0.926"
$ws.Range("E50").Value = "  -2.05%  "

$ws.Range("D51").Value = "8.93"
$ws.Range("E51").Value = "  +0.12%  "
